$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 180849.53
$ws.Range("J17").Value = 184107.47
$ws.Range("L17").Value = 552322.41
$ws.Range("N17").Value = -552658.41
$ws.Range("H133").Value = 39610.168
$ws.Range("J133").Value = 39610.168
$ws.Range("L133").Value = 39610.168
$ws.Range("N133").Value = -49730.168
$ws.Range("H135").Value = 2036.2916
$ws.Range("I135").Value = 823.2778
$ws.Range("K135").Value = 7409.500199999999
$ws.Range("M135").Value = -4874.500199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1271.6765
$ws.Range("I2").Value = 1270.4231
$ws.Range("K2").Value = 1270.4231
$ws.Range("M2").Value = -1157.4231
$ws.Range("H32").Value = 9802.718000000001
$ws.Range("I32").Value = 10210.024
$ws.Range("J32").Value = 6462.8
$ws.Range("K32").Value = 10210.024
$ws.Range("L32").Value = 6462.8
$ws.Range("M32").Value = -9923.023999999999
$ws.Range("N32").Value = -7036.8
$ws.Range("H116").Value = 1271.6765
$ws.Range("I116").Value = 1270.4231
$ws.Range("K116").Value = 1270.4231
$ws.Range("M116").Value = 1023.5769

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1271.6765
$ws.Range("I3").Value = 1270.4231
$ws.Range("K3").Value = 1270.4231
$ws.Range("M3").Value = -1156.4231
$ws.Range("H99").Value = 2881.611
$ws.Range("I99").Value = 2134.9285
$ws.Range("J99").Value = 5495
$ws.Range("K99").Value = 2134.9285
$ws.Range("L99").Value = 5495
$ws.Range("M99").Value = -636.9285
$ws.Range("N99").Value = -8491
$ws.Range("H105").Value = 3499.6562
$ws.Range("I105").Value = 3417.9092
$ws.Range("J105").Value = 3679.5
$ws.Range("K105").Value = 3417.9092
$ws.Range("L105").Value = 3679.5
$ws.Range("M105").Value = -1670.9092
$ws.Range("N105").Value = -7173.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1872.6578
$ws.Range("I31").Value = 1044.931
$ws.Range("K31").Value = 1044.931
$ws.Range("M31").Value = -749.931
$ws.Range("H34").Value = 1872.6578
$ws.Range("I34").Value = 1044.931
$ws.Range("K34").Value = 1044.931
$ws.Range("M34").Value = -842.931
$ws.Range("H134").Value = 4291.381
$ws.Range("I134").Value = 1920
$ws.Range("J134").Value = 6899.9
$ws.Range("K134").Value = 5760
$ws.Range("L134").Value = 20699.7
$ws.Range("M134").Value = -3225
$ws.Range("N134").Value = -25769.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I32").Value = 50000250
$ws.Range("J32").Value = 2050.3333
$ws.Range("K32").Value = 150000750
$ws.Range("L32").Value = 6150.999899999999
$ws.Range("M32").Value = -150000467
$ws.Range("N32").Value = -6716.999899999999
$ws.Range("H46").Value = 284854.53
$ws.Range("I46").Value = 208480
$ws.Range("J46").Value = 348500
$ws.Range("K46").Value = 625440
$ws.Range("L46").Value = 1045500
$ws.Range("M46").Value = -625349
$ws.Range("N46").Value = -1045682
$ws.Range("H44").Value = 999999
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 999999
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 2999997
$ws.Range("N44").Value = -3000793
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("M48").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4244.517
$ws.Range("I80").Value = 4269.1113
$ws.Range("K80").Value = 4269.1113
$ws.Range("M80").Value = -3271.1113
$ws.Range("H83").Value = 4244.517
$ws.Range("I83").Value = 4269.1113
$ws.Range("K83").Value = 21345.5565
$ws.Range("M83").Value = -16353.5565
$ws.Range("H104").Value = 63662.332
$ws.Range("J104").Value = 63662.332
$ws.Range("L104").Value = 63662.332
$ws.Range("N104").Value = -70650.33199999999
$ws.Range("H107").Value = 923.3077
$ws.Range("J107").Value = 1294
$ws.Range("L107").Value = 1294
$ws.Range("N107").Value = -5134
$ws.Range("H132").Value = 4632096.5
$ws.Range("I132").Value = 5749494.5
$ws.Range("J132").Value = 2876.5715
$ws.Range("K132").Value = 17248483.5
$ws.Range("L132").Value = 8629.7145
$ws.Range("M132").Value = -17245953.5
$ws.Range("N132").Value = -13689.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1814.6296
$ws.Range("I22").Value = 317.44446
$ws.Range("J22").Value = 2563.2222
$ws.Range("K22").Value = 317.44446
$ws.Range("L22").Value = 2563.2222
$ws.Range("M22").Value = -22.44445999999999
$ws.Range("N22").Value = -3153.2222
$ws.Range("H27").Value = 1814.6296
$ws.Range("I27").Value = 317.44446
$ws.Range("J27").Value = 2563.2222
$ws.Range("K27").Value = 317.44446
$ws.Range("L27").Value = 2563.2222
$ws.Range("M27").Value = -210.44446
$ws.Range("N27").Value = -2777.2222
$ws.Range("H128").Value = 64996.5
$ws.Range("J128").Value = 64996.5
$ws.Range("L128").Value = 64996.5
$ws.Range("N128").Value = -74956.5
$ws.Range("H136").Value = 3864.4827
$ws.Range("I136").Value = 3133.8333
$ws.Range("J136").Value = 5060.091
$ws.Range("K136").Value = 9401.499899999999
$ws.Range("L136").Value = 15180.273
$ws.Range("M136").Value = -6851.499899999999
$ws.Range("N136").Value = -20280.273

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3649.1428
$ws.Range("J96").Value = 3924
$ws.Range("L96").Value = 3924
$ws.Range("N96").Value = -6670
$ws.Range("H100").Value = 1980.2
$ws.Range("I100").Value = 850.25
$ws.Range("K100").Value = 1700.5
$ws.Range("M100").Value = -1159.5
$ws.Range("H107").Value = 768.85
$ws.Range("I107").Value = 446.5
$ws.Range("J107").Value = 1091.2
$ws.Range("K107").Value = 1339.5
$ws.Range("L107").Value = 3273.6
$ws.Range("M107").Value = 580.5
$ws.Range("N107").Value = -7113.6
$ws.Range("H124").Value = 111809.664
$ws.Range("J124").Value = 111809.664
$ws.Range("L124").Value = 111809.664
$ws.Range("N124").Value = -121629.664
$ws.Range("H133").Value = 106499.5
$ws.Range("J133").Value = 106499.5
$ws.Range("L133").Value = 106499.5
$ws.Range("N133").Value = -116619.5
$ws.Range("H136").Value = 6478.8423
$ws.Range("I136").Value = 9184
$ws.Range("J136").Value = 1276.6154
$ws.Range("K136").Value = 27552
$ws.Range("L136").Value = 9184
$ws.Range("M136").Value = -25002
$ws.Range("N136").Value = -8929.8462
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
